$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2 to 53 with new date (col A) and value (col B) data
$data = @(
    @(39400, 3.582807286231798),
    @(39583, -3.395106996021084),
    @(39765, 3.499081613727355),
    @(39948, 2.643706631981502),
    @(40130, 2.323758253012315),
    @(40310, 1.86889045152472),
    @(40494, 1.711273044056469),
    @(40676, 1.792365984063807),
    @(40862, 1.894327589997687),
    @(41044, -1.203281951233052),
    @(41228, 1.349970205759888),
    @(41409, -1.129543206172372),
    @(41592, 1.457269272787158),
    @(41774, -0.745755750114057),
    @(41957, 0.06404262623020429),
    @(42137, 0.2463991311210521),
    @(42321, 0.6724301216571575),
    @(42503, 0.6559243910605232),
    @(42689, 0.5936072881632271),
    @(42867, 0.3249989166702818),
    @(43053, 0.8673551286853183),
    @(43145, 0.6331942894404392),
    @(43235, 1.523252944018139),
    @(43326, -0.5201366209837346),
    @(43418, 0.3),
    @(43510, 0.1895219545376108),
    @(43600, -1.1),
    @(43691, -0.5586173343648539),
    @(43783, 0.9418773066947779),
    @(43875, 1.001530035891491),
    @(43966, -0.3),
    @(44068, 0.5863322451846074),
    @(44159, 1.879044851730669),
    @(44251, 0.5891195155317774),
    @(44341, -1.620076036519961),
    @(44432, 1.285713798234809),
    @(44525, 0.9260365910423758),
    @(44617, 1.323876910632023),
    @(44706, -2.136777872354216),
    @(44798, 0.7934051314956321),
    @(44890, 1.029975481959951),
    @(44981, 0.7043810499788776),
    @(45071, 0.1406555772994125),
    @(45163, 0.4268111723034735),
    @(45254, 0.3049656067824742),
    @(45345, 0.7386451510207621),
    @(45436, -1.147231087866345),
    @(45534, 0.5563279776301044),
    @(45618, 0.6119095091543301),
    @(45713, -0.1750403534197602),
    @(45800, 2.014657057377804),
    @(45891, 1.176666004305858),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $data[$i][1]
}

# Remove old rows 54 to 73 (no longer present in the new data range)
$ws.Range("A54:B73").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
